$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '43.696.54'
$ws.Range('E2').Value = '  -0.16%  '

# Row 3
$ws.Range('D3').Value = '2.289.24'
$ws.Range('E3').Value = '  -1.03%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
Set-TextValue $ws.Range('D5') '103.52'
$ws.Range('E5').Value = '  +6.92%  '

# Row 6
Set-TextValue $ws.Range('D6') '270.72'
$ws.Range('E6').Value = '  -0.32%  '

# Row 7
$ws.Range('E7').Value = '  -0.45%  '

# Row 8
$ws.Range('E8').Value = '  -0.09%  '

# Row 9
$ws.Range('E9').Value = '  -3.00%  '

# Row 10
Set-TextValue $ws.Range('D10') '45.88'
$ws.Range('E10').Value = '  +1.21%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0929'
$ws.Range('E11').Value = '  -1.91%  '

# Row 12
Set-TextValue $ws.Range('D12') '7.98'
$ws.Range('E12').Value = '  -0.96%  '

# Row 13
$ws.Range('E13').Value = '  +1.62%  '

# Row 14
Set-TextValue $ws.Range('D14') '15.58'
$ws.Range('E14').Value = '  +0.35%  '

# Row 15
$ws.Range('E15').Value = '  -2.19%  '

# Row 16
$ws.Range('D16').Value = '2.287.76'
$ws.Range('E16').Value = '  -1.25%  '

# Row 17
$ws.Range('D17').Value = '43.672.01'
$ws.Range('E17').Value = '  -0.16%  '

# Row 18
$ws.Range('E18').Value = '  +0.83%  '

# Row 19
$ws.Range('E19').Value = '  -2.11%  '

# Row 20
Set-TextValue $ws.Range('D20') '72.22'
$ws.Range('E20').Value = '  -1.48%  '

# Row 21
$ws.Range('E21').Value = '  +9.89%  '

# Row 22
Set-TextValue $ws.Range('D22') '231.82'
$ws.Range('E22').Value = '  -3.35%  '

# Row 23
Set-TextValue $ws.Range('D23') '2.89'
$ws.Range('E23').Value = '  +13.84%  '

# Row 24
Set-TextValue $ws.Range('D24') '9.13'
$ws.Range('E24').Value = '  -3.00%  '

# Row 25
$ws.Range('E25').Value = '  -0.05%  '

# Row 26
Set-TextValue $ws.Range('D26') '11.18'
$ws.Range('E26').Value = '  -1.70%  '

# Row 27
$ws.Range('E27').Value = '  -0.81%  '

# Row 28
Set-TextValue $ws.Range('D28') '40.38'
$ws.Range('E28').Value = '  +6.79%  '

# Row 29
$ws.Range('E29').Value = '  -2.56%  '

# Row 30
Set-TextValue $ws.Range('D30') '177.11'
$ws.Range('E30').Value = '  +1.55%  '

# Row 31
Set-TextValue $ws.Range('D31') '21.78'
$ws.Range('E31').Value = '  -2.84%  '

# Row 32
$ws.Range('E32').Value = '  -0.68%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D33') '5.46'
$ws.Range('E33').Value = '  -0.10%  '

# Row 34
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D34') '4.91'
$ws.Range('E34').Value = '  +11.91%  '

# Row 35
$ws.Range('E35').Value = '  -0.23%  '

# Row 36
$ws.Range('E36').Value = '  +0.02%  '

# Row 37
$ws.Range('E37').Value = '  -2.54%  '

# Row 38
Set-TextValue $ws.Range('D38') '3.58'
$ws.Range('E38').Value = '  +6.26%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.237'
$ws.Range('E39').Value = '  -3.63%  '

# Row 40
Set-TextValue $ws.Range('D40') '2.32'
$ws.Range('E40').Value = '  -0.86%  '

# Row 41
Set-TextValue $ws.Range('D41') '1.37'
$ws.Range('E41').Value = '  +0.19%  '

# Row 42
Set-TextValue $ws.Range('D42') '12.32'
$ws.Range('E42').Value = '  +1.23%  '

# Row 43
Set-TextValue $ws.Range('D43') '65.66'
$ws.Range('E43').Value = '  +5.04%  '

# Row 44
Set-TextValue $ws.Range('D44') '5.25'
$ws.Range('E44').Value = '  -1.92%  '

# Row 45
Set-TextValue $ws.Range('D45') '8.77'
$ws.Range('E45').Value = '  -4.47%  '

# Row 46
$ws.Range('E46').Value = '  -1.15%  '

# Row 47
$ws.Range('E47').Value = '  +2.30%  '

# Row 48
$ws.Range('E48').Value = '  -1.17%  '

# Row 49
Set-TextValue $ws.Range('D49') '0.447'
$ws.Range('E49').Value = '  +8.20%  '

# Row 50
Set-TextValue $ws.Range('D50') '1.53'
$ws.Range('E50').Value = '  +10.75%  '

# Row 51
$ws.Range('D51').Value = '2.512.71'
$ws.Range('E51').Value = '  -1.08%  '
